$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in column E per diff
$ws.Range("E2").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("E5").Value = 72

# E6 value changes and loses its number-format style (s="1" -> default)
$ws.Range("E6").Value = 72
$ws.Range("E6").ClearFormats()

$ws.Range("E7").Value = 144
$ws.Range("E8").Value = 72
$ws.Range("E9").Value = 72
$ws.Range("E10").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("E13").Value = 8640
$ws.Range("E14").Value = 8640
$ws.Range("E15").Value = 17280
$ws.Range("E16").Value = 8640
$ws.Range("E17").Value = 8640

# Select entire column F (producing <selection activeCell="F1" sqref="F1:F1048576"/>)
$ws.Columns("F").Select()
